$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.092.96"
$ws.Range("E2").Value = "  -3.60%  "

$ws.Range("D3").Value = "3.385.83"
$ws.Range("E3").Value = "  -4.46%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'563.48"
$ws.Range("E5").Value = "  -3.85%  "

$ws.Range("D6").Value = "'185.80"
$ws.Range("E6").Value = "  -6.03%  "

$ws.Range("E7").Value = "  -2.14%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "3.379.11"
$ws.Range("E9").Value = "  -4.16%  "

$ws.Range("D10").Value = "'0.191"
$ws.Range("E10").Value = "  -7.92%  "

$ws.Range("D11").Value = "'0.600"
$ws.Range("E11").Value = "  -4.64%  "

$ws.Range("D12").Value = "'48.44"
$ws.Range("E12").Value = "  -7.12%  "

$ws.Range("D13").Value = "'0.0000274"
$ws.Range("E13").Value = "  -5.64%  "

$ws.Range("D14").Value = "'8.83"
$ws.Range("E14").Value = "  -5.58%  "

$ws.Range("D15").Value = "3.926.25"
$ws.Range("E15").Value = "  -4.31%  "

$ws.Range("D16").Value = "'608.95"
$ws.Range("E16").Value = "  -10.49%  "

$ws.Range("D17").Value = "66.944.60"
$ws.Range("E17").Value = "  -3.88%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'18.14"
$ws.Range("E18").Value = "  -2.58%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.377.25"
$ws.Range("E19").Value = "  -4.30%  "

$ws.Range("E20").Value = "  -2.94%  "

$ws.Range("D21").Value = "'11.75"
$ws.Range("E21").Value = "  -5.70%  "

$ws.Range("D22").Value = "'0.922"
$ws.Range("E22").Value = "  -5.06%  "

$ws.Range("D23").Value = "'17.12"
$ws.Range("E23").Value = "  -4.57%  "

$ws.Range("D24").Value = "'5.17"
$ws.Range("E24").Value = "  -1.70%  "

$ws.Range("D25").Value = "'99.11"
$ws.Range("E25").Value = "  -8.38%  "

$ws.Range("D26").Value = "'4.13"
$ws.Range("E26").Value = "  -6.44%  "

$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("D28").Value = "'2.76"
$ws.Range("E28").Value = "  -6.42%  "

$ws.Range("D29").Value = "'9.59"
$ws.Range("E29").Value = "  -7.31%  "

$ws.Range("D30").Value = "'8.91"
$ws.Range("E30").Value = "  -8.31%  "

$ws.Range("D31").Value = "'31.05"
$ws.Range("E31").Value = "  -7.28%  "

$ws.Range("D32").Value = "'4.00"
$ws.Range("E32").Value = "  -8.61%  "

$ws.Range("D33").Value = "'6.40"
$ws.Range("E33").Value = "  -7.77%  "

$ws.Range("D34").Value = "'11.26"
$ws.Range("E34").Value = "  -5.84%  "

$ws.Range("D35").Value = "'555.35"
$ws.Range("E35").Value = "  +10.20%  "

$ws.Range("E36").Value = "  -4.71%  "

$ws.Range("D37").Value = "3.888.18"
$ws.Range("E37").Value = "  +2.13%  "

$ws.Range("D38").Value = "'58.74"
$ws.Range("E38").Value = "  -5.92%  "

$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("D40").Value = "'3.47"
$ws.Range("E40").Value = "  -5.42%  "

$ws.Range("D41").Value = "'3.57"
$ws.Range("E41").Value = "  +32.81%  "

$ws.Range("D42").Value = "0.0₃0733"
$ws.Range("E42").Value = "  -10.96%  "

$ws.Range("E43").Value = "  -7.33%  "

$ws.Range("D44").Value = "'0.130"
$ws.Range("E44").Value = "  -5.10%  "

$ws.Range("D45").Value = "'0.353"
$ws.Range("E45").Value = "  -5.46%  "

$ws.Range("D46").Value = "'32.67"
$ws.Range("E46").Value = "  -6.71%  "

$ws.Range("D47").Value = "'0.0424"
$ws.Range("E47").Value = "  -8.09%  "

$ws.Range("D48").Value = "'3.27"
$ws.Range("E48").Value = "  -3.01%  "

$ws.Range("D49").Value = "'2.71"
$ws.Range("E49").Value = "  -8.54%  "

$ws.Range("D50").Value = "'0.132"
$ws.Range("E50").Value = "  -4.42%  "

$ws.Range("E51").Value = "  -0.21%  "
